$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("Old") updates
$ws.Range("B2").Value = 64
$ws.Range("Q2").Value = 53
$ws.Range("R2").Value = 147
$ws.Range("S2").Value = 118
$ws.Range("T2").Value = 132
$ws.Range("U2").Value = 52

# Row 3 ("New") updates
$ws.Range("B3").Value = 76
$ws.Range("Q3").Value = 102
$ws.Range("R3").Value = 235
$ws.Range("S3").Value = 177
$ws.Range("T3").Value = 189
$ws.Range("U3").Value = 67
